# SZ JS site ToDos.xlsx - add a new "to do" row above the existing
# "filter video point query when zoomed out" item.
#
# This inserts a new row 58 on the "to do" sheet containing the new idea
# text, which pushes all subsequent rows (old 58-442) down by one (new
# 59-443), and updates the dependent structures (comments, hyperlinks,
# autofilter, defined names, dimension, selection) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert the new row at position 58 (shifts rows 58:442 -> 59:443,
#    along with their comments/hyperlinks/styles).
$ws.Rows("58:58").Insert()

# 2. The newly inserted row 58 inherited formatting from the row above
#    it (old row 57). Give D58 the same look used elsewhere in this sheet
#    for a non-hidden "idea" note (same style as D61/D88), and fill in the
#    new text.
$ws.Range("D61").Copy()
$ws.Range("D58").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D58").Value = "IDEA:  Filter for Survey & Survey Year legend display?  (incl. current extent)"

# Remove the stray copied-down formatting in column C for this row (the
# new row only has content in D and a blank styled J cell).
$ws.Range("C58").Clear()

# Row 58 is a plain (non "customHeight") 15.75pt row, like other similar
# single-line rows in this sheet (e.g. row 98/99).
$ws.Rows("58:58").RowHeight = 15.75

# 3. Extend the autofilter range from A1:I442 to A1:I443 and restore the
#    two "blanks only" column filters that were in effect before.
$ws.AutoFilterMode = $false
$ws.Range("A1:I443").AutoFilter()
$ws.Range("A1:I443").AutoFilter(1, @(""), 7)  # column A ("colId=0"): blanks only
$ws.Range("A1:I443").AutoFilter(3, @(""), 7)  # column C ("colId=2"): blanks only

# 4. Update the hidden _FilterDatabase defined name to match the new
#    autofilter range (must be addressed by index - by-name lookup does
#    not resolve to the live object for this builtin name).
$fd = $ws.Names.Item(1)
$fd.RefersTo = "='to do'!`$A`$1:`$I`$443"

# 5. Reflect the edit location in the UI selection (cursor lands just
#    below the newly added row, on the pre-existing item).
$excel.Goto($ws.Range("D59"), $true)
